# 25-11-2022: Aalborg EnagyPLAN 15.1 path fixed, ESD calculation bug fixed
#
# This script:
#  1. Adds a new "Sheet2" after "Sheet1", populated with the previous
#     (pre-edit) model layout but re-run with a different set of input
#     assumptions (plus an extra helper row 20).
#  2. Updates "Sheet1" in place: new input values, B8 becomes a SUM()
#     over the capacity-factor table, a new "Offshore"/"pv" pair of rows
#     is added to that table (rows 11-12), and the weighted-PEF helper
#     (previously G12/H12) is relocated down to a new row 15 so it can
#     pick up the two new rows.
#  3. Restores Sheet1 as the active/selected sheet and leaves the B14
#     cell selected on both sheets, matching the authored workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create Sheet2 (placed right after Sheet1) and fill it in with the
#    model as it used to read on Sheet1 before this commit, but using
#    the new set of input numbers for this sheet.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("B1").Value = "TWh"

$ws2.Range("A2").Value = "el_i"
$ws2.Range("B2").Value = 0.27

$ws2.Range("A3").Value = "PEF_i^el"
$ws2.Range("B3").Formula = "=1/0.47"

$ws2.Range("A4").Value = "Ngas_i"
$ws2.Range("B4").Value = 0.04

$ws2.Range("A5").Value = "BM_i"
$ws2.Range("B5").Value = 0

$ws2.Range("A6").Value = "WS_i"
$ws2.Range("B6").Value = 0
$ws2.Range("G6").Value = "Local PEF caculation"
$ws2.Range("G6").Font.Bold = $true

$ws2.Range("G7").Value = "production (TWh)"
$ws2.Range("H7").Value = "eff"
$ws2.Range("I7").Value = "PEF"

$ws2.Range("A8").Value = "el_lp"
$ws2.Range("B8").Value = 0.81
$ws2.Range("F8").Value = "wind "
$ws2.Range("G8").Value = 0.67
$ws2.Range("I8").Value = 1

$ws2.Range("A9").Value = "el_e"
$ws2.Range("B9").Value = 1.07
$ws2.Range("F9").Value = "chp"
$ws2.Range("G9").Value = 0.04
$ws2.Range("H9").Value = 0.25
$ws2.Range("I9").Formula = "=1/H9"

$ws2.Range("A10").Value = "PEF _lp^el"
$ws2.Range("B10").Formula = "=H12"
$ws2.Range("F10").Value = "cshp"
$ws2.Range("G10").Value = 0.1
$ws2.Range("H10").Value = 0.54
$ws2.Range("I10").Formula = "=1/H10"

$ws2.Range("A11").Value = "BM_c"
$ws2.Range("B11").Value = 1.05

$ws2.Range("A12").Value = "WS_c"
$ws2.Range("B12").Value = 0.5
$ws2.Range("G12").Value = "PEF _lp^el"
$ws2.Range("H12").Formula = "=(G8*I8+G9*I9+G10*I10)/B8"

$ws2.Range("A13").Value = "H_HP"
$ws2.Range("B13").Value = 0.5

$ws2.Range("A14").Value = "COP"
$ws2.Range("B14").Value = 3.6

$ws2.Range("A18").Value = "ESD"
$ws2.Range("B18").Formula = "=(B4+B5+B6)/((B8-B9)*B10+B4+B11+B12+(B13-(B13/B14)))"

$ws2.Range("G20").Formula = "=G8*1 + G9 * (1/0.25) + G10 * (1/0.54)"
$ws2.Range("H20").Formula = "=1/0.25"

# ---------------------------------------------------------------------
# 2. Update Sheet1 with the new inputs / restructured PEF table.
# ---------------------------------------------------------------------

# el_i unchanged, Ngas_i -> 0.51
$ws1.Range("B4").Value = 0.51

# el_lp: capacity-factor table now spans G8:G12 (Offshore + pv rows
# added below), so the formula is rewritten as a SUM() over it.
$ws1.Range("G8").Value = 0.57
$ws1.Range("B8").Formula = "=SUM(G8:G12)"

# el_e
$ws1.Range("B9").Value = 0.29
$ws1.Range("G9").Value = 0.3

# PEF _lp^el now looks up the relocated helper row (H15 instead of H12)
$ws1.Range("B10").Formula = "=H15"

# BM_c, plus new "pv" row feeding the capacity-factor table
$ws1.Range("B11").Value = 1.52
$ws1.Range("F11").Value = "pv"
$ws1.Range("G11").Value = 0.53
$ws1.Range("I11").Value = 1

# WS_c, plus new "Offshore" row feeding the capacity-factor table.
# The old weighted-PEF helper that used to live in G12/H12 is relocated
# to row 15 below, so clear H12's old formula first.
$ws1.Range("H12").ClearContents()
$ws1.Range("B12").Value = 0.5
$ws1.Range("F12").Value = "Offshore"
$ws1.Range("G12").Value = 0.02
$ws1.Range("I12").Value = 1

# H_HP
$ws1.Range("B13").Value = 0.33

# Relocate the weighted-PEF helper (used to live at G12/H12) down to a
# fresh row 15 so it can include the two new capacity-factor rows.
$ws1.Range("G15").Value = "PEF _lp^el"
$ws1.Range("H15").Formula = "=(G8*I8+G9*I9+G10*I10+G11*I11+G12*I12)/B8"

# ---------------------------------------------------------------------
# 3. Restore selections / active sheet so the saved view matches.
# ---------------------------------------------------------------------
$ws2.Range("B14").Select()
$ws1.Activate()
$ws1.Range("B14").Select()
